$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2..3) {
    $ws.Range("AG$row").Value = 1094172
    $ws.Range("AR$row").Value = 382
    $ws.Range("AS$row").Value = 1095239
    $ws.Range("AU$row").Value = 624000
    $ws.Range("AV$row").Value = 624000
    $ws.Range("BI$row").Value = "xq03"
    $ws.Range("BJ$row").Value = "Chụp Xquang bể thận-niệu quản xuôi dòng [Số hóa]"
    $ws.Range("BX$row").Value = 124800
    $ws.Range("BZ$row").Value = 499200
}
